$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet marks completed tasks with a Wingdings "ü" glyph, which renders as a
# check mark. Some rows use a plain (no-fill) check style, others a gray-fill
# check style (matching the gray-fill "not applicable" columns used elsewhere).
$chk = [char]252

# --- Row 19 ---
# G19: blank -> checkmark, no fill (same look as G3/G4/G5 etc.)
$ws.Range("G19").Font.Name = "Wingdings"
$ws.Range("G19").Value = $chk

# H19: blank -> checkmark, gray fill (same look as F19/H3 etc.)
$ws.Range("H19").Font.Name = "Wingdings"
$ws.Range("H19").Interior.Pattern = -4124       # xlPatternGray (matches existing gray-fill cells)
$ws.Range("H19").Interior.Color = 11579568      # matches the existing gray fill (RGB 00B0B0B0)
$ws.Range("H19").Value = $chk

# --- Row 24 ---
# G24: blank -> checkmark, no fill (same look as G3/G4/G5 etc.)
$ws.Range("G24").Font.Name = "Wingdings"
$ws.Range("G24").Value = $chk

# --- Row 27 ---
# E27: blank -> checkmark, no fill (same look as E3 etc.)
$ws.Range("E27").Font.Name = "Wingdings"
$ws.Range("E27").Value = $chk

# --- Row 34 ---
# H34: blank -> checkmark, gray fill (same look as F34/H3 etc.)
$ws.Range("H34").Font.Name = "Wingdings"
$ws.Range("H34").Interior.Pattern = -4124       # xlPatternGray (matches existing gray-fill cells)
$ws.Range("H34").Interior.Color = 11579568      # matches the existing gray fill (RGB 00B0B0B0)
$ws.Range("H34").Value = $chk
